$wb = $excel.ActiveWorkbook

# --- Sheet "Variables": insert a new row 2 for the "ident" variable ---
$wsVar = $wb.Worksheets.Item("Variables")
$wsVar.Select()
$wsVar.Rows.Item(2).Insert()
$wsVar.Range("B2").Value = "ident"
$wsVar.Range("C2").Value = "participant identifier"
$wsVar.Range("D2").Value = "integer"

# restore the selection on the Variables sheet to H9
$wsVar.Range("H9").Select()

# --- Sheet "Categories": move selection / viewport ---
$wsCat = $wb.Worksheets.Item("Categories")
$wsCat.Select()
$excel.ActiveWindow.ScrollRow = 19
$wsCat.Range("A2:XFD2").Select()
